$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.052.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.90%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.687.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.52%  '

$ws.Range("E6").Value = '  +2.30%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.33'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.267'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0639'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.929.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.684.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.76%  '

$ws.Range("E15").Value = '  +3.40%  '

$ws.Range("E16").Value = '  +7.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.028.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.67%  '

$ws.Range("E20").Value = '  +2.40%  '

$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.66%  '

$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("E25").Value = '  +0.84%  '

$ws.Range("E26").Value = '  +2.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.113'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.13%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("E30").Value = '  +2.08%  '

$ws.Range("E31").Value = '  +4.04%  '

$ws.Range("E32").Value = '  +4.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.514.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.37%  '

$ws.Range("E35").Value = '  +2.70%  '

$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '83.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.614'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0180'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.30%  '

$ws.Range("E40").Value = '  -3.55%  '

$ws.Range("E41").Value = '  +0.63%  '

$ws.Range("E42").Value = '  +2.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.843'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.93%  '

$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("E45").Value = '  +3.27%  '

$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.17%  '

$ws.Range("E48").Value = '  +4.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.817.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.39%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0120'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.21%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '93.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.88%  '

